# This script "refreshes" the timestamps baked into the e-commerce operations
# tracking workbook, moving every captured "Generated"/"last updated" style
# timestamp forward from 12:14 (xx:14) to 12:18 (xx:18) on 2025-10-06 - as if
# the workbook had been regenerated ~4 minutes later.

$wb = $excel.ActiveWorkbook

# --- Dashboard sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("Dashboard")

# Header "Generated: ..." timestamp (includes seconds).
$ws.Range("A2").Value = "Generated: 2025-10-06 12:18:05"

# "Last Updated" column (E) for each KPI row.
$ws.Range("E6").Value  = "2025-10-06 12:18"
$ws.Range("E7").Value  = "2025-10-06 12:18"
$ws.Range("E8").Value  = "2025-10-06 12:18"
$ws.Range("E9").Value  = "2025-10-06 12:18"
$ws.Range("E10").Value = "2025-10-06 12:18"
$ws.Range("E11").Value = "2025-10-06 12:18"
$ws.Range("E12").Value = "2025-10-06 12:18"

# --- Bash Queries Response sheet -------------------------------------------
$ws = $wb.Worksheets.Item("Bash Queries Response")

$ws.Range("B3").Value = "2025-10-06 12:18"
$ws.Range("B4").Value = "2025-10-06 12:18"
$ws.Range("B5").Value = "2025-10-06 12:18"

# --- Stock Replenishment sheet ----------------------------------------------
$ws = $wb.Worksheets.Item("Stock Replenishment")

$ws.Range("G3").Value = "2025-10-06 09:18"
$ws.Range("H3").Value = "2025-10-06 11:18"
$ws.Range("G4").Value = "2025-10-06 07:18"
$ws.Range("H4").Value = "2025-10-06 12:18"
$ws.Range("G5").Value = "2025-10-06 10:18"

# --- System Errors sheet ----------------------------------------------------
$ws = $wb.Worksheets.Item("System Errors")

$ws.Range("B3").Value = "2025-10-06 12:18"
$ws.Range("B4").Value = "2025-10-06 12:18"
$ws.Range("B5").Value = "2025-10-06 12:18"
